# Populate the "Lương" (Salary) worksheet with the per-employee salary
# breakdown. The sheet is the second sheet in the workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

$rows = @(
    @("Danh mục", 8),
    @("Ngày công", 11),
    @("Phụ cấp", 385000),
    @("Lương cơ bản tại CẦN THƠ", 0),
    @("Chiết khấu sale chính tại CẦN THƠ", 0),
    @("Chiết khấu sale phụ tại CẦN THƠ", 0),
    @("Đơn 1 bác sĩ tại CẦN THƠ", 0),
    @("Đơn 2 bác sĩ tại CẦN THƠ", 0),
    @("Công phụ phẫu 1 tại CẦN THƠ", 0),
    @("Công phụ phẫu 2 tại CẦN THƠ", 0),
    @("Lương cơ bản tại LONG XUYÊN", 0),
    @("Chiết khấu sale chính tại LONG XUYÊN", 0),
    @("Chiết khấu sale phụ tại LONG XUYÊN", 0),
    @("Đơn 1 bác sĩ tại LONG XUYÊN", 0),
    @("Đơn 2 bác sĩ tại LONG XUYÊN", 0),
    @("Công phụ phẫu 1 tại LONG XUYÊN", 0),
    @("Công phụ phẫu 2 tại LONG XUYÊN", 0),
    @("Lương cơ bản tại SÓC TRĂNG", 0),
    @("Chiết khấu sale chính tại SÓC TRĂNG", 0),
    @("Chiết khấu sale phụ tại SÓC TRĂNG", 0),
    @("Đơn 1 bác sĩ tại SÓC TRĂNG", 0),
    @("Đơn 2 bác sĩ tại SÓC TRĂNG", 0),
    @("Công phụ phẫu 1 tại SÓC TRĂNG", 0),
    @("Công phụ phẫu 2 tại SÓC TRĂNG", 0)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 1
    $label = $rows[$i][0]
    $amount = $rows[$i][1]
    $ws.Cells.Item($r, 1).Value = $label
    $ws.Cells.Item($r, 2).Value = $amount
}
